$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.007.11"
$ws.Range("E2").Value = "  +5.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.705.98"
$ws.Range("E3").Value = "  +6.79%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.84"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.76"
$ws.Range("E6").Value = "  +1.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.700.10"
$ws.Range("E7").Value = "  +6.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.644"
$ws.Range("E8").Value = "  +1.12%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.765"
$ws.Range("E10").Value = "  -2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.181"
$ws.Range("E11").Value = "  +10.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000392"
$ws.Range("E12").Value = "  +46.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.96"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.36"
$ws.Range("E14").Value = "  +5.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.306.77"
$ws.Range("E15").Value = "  +7.20%  "

$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.59"
$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.713.63"
$ws.Range("E18").Value = "  +6.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.05"
$ws.Range("E19").Value = "  +5.26%  "

$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.108.17"
$ws.Range("E21").Value = "  +5.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "448.59"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.98"
$ws.Range("E23").Value = "  +18.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.65"
$ws.Range("E24").Value = "  -1.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("E25").Value = "  -3.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.01"
$ws.Range("E26").Value = "  +11.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.32"
$ws.Range("E27").Value = "  +1.26%  "

$ws.Range("E28").Value = "  +0.64%  "

$ws.Range("E29").Value = "  +4.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.57"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +7.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.69"
$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.32"
$ws.Range("E33").Value = "  -3.68%  "

$ws.Range("E34").Value = "  -1.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.18"
$ws.Range("E35").Value = "  +3.00%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  -1.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0491"
$ws.Range("E38").Value = "  -2.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("E39").Value = "  +40.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0750"
$ws.Range("E40").Value = "  +12.88%  "

$ws.Range("E41").Value = "  +4.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "29.11"
$ws.Range("E42").Value = "  +32.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.41"
$ws.Range("E44").Value = "  +1.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.16"
$ws.Range("E45").Value = "  +2.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.10"
$ws.Range("E46").Value = "  +4.79%  "

$ws.Range("E47").Value = "  -5.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.39"
$ws.Range("E48").Value = "  -2.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.63"
$ws.Range("E49").Value = "  -5.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.307"
$ws.Range("E50").Value = "  -3.77%  "

$ws.Range("E51").Value = "  +16.68%  "
